$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 1002
$ws.Range("I74").Value = 1002
$ws.Range("K74").Value = 1002
$ws.Range("M74").Value = -66

# Row 77
$ws.Range("H77").Value = 1002
$ws.Range("I77").Value = 1002
$ws.Range("K77").Value = 5010
$ws.Range("M77").Value = -330

# Row 86
$ws.Range("H86").Value = 7491.2144
$ws.Range("I86").Value = 12799.5
$ws.Range("J86").Value = 5367.9
$ws.Range("K86").Value = 12799.5
$ws.Range("L86").Value = 5367.9
$ws.Range("M86").Value = -11676.5
$ws.Range("N86").Value = -7613.9

# Row 89
$ws.Range("H89").Value = 7491.2144
$ws.Range("I89").Value = 12799.5
$ws.Range("J89").Value = 5367.9
$ws.Range("K89").Value = 63997.5
$ws.Range("L89").Value = 26839.5
$ws.Range("M89").Value = -58381.5
$ws.Range("N89").Value = -38071.5

# Row 103
$ws.Range("H103").Value = 1444.8334
$ws.Range("I103").Value = 1157.2222
$ws.Range("J103").Value = 1732.4445
$ws.Range("K103").Value = 3471.6666
$ws.Range("L103").Value = 5197.333500000001
$ws.Range("M103").Value = -2885.6666
$ws.Range("N103").Value = -6369.333500000001

# Row 106
$ws.Range("H106").Value = 3448.8333
$ws.Range("I106").Value = 2944.182
$ws.Range("K106").Value = 2944.182
$ws.Range("M106").Value = -2313.182

# Row 112
$ws.Range("H112").Value = 1940.1875
$ws.Range("J112").Value = 2057.2307
$ws.Range("L112").Value = 6171.6921
$ws.Range("N112").Value = -8387.6921

# Row 132
$ws.Range("H132").Value = 4865.279
$ws.Range("I132").Value = 2431.6155
$ws.Range("K132").Value = 7294.8465
$ws.Range("M132").Value = -4764.8465

# Row 137
$ws.Range("H137").Value = 5458.9585
$ws.Range("I137").Value = 5668.143
$ws.Range("J137").Value = 3994.6667
$ws.Range("K137").Value = 17004.429
$ws.Range("L137").Value = 11984.0001
$ws.Range("M137").Value = -14454.429
$ws.Range("N137").Value = -17084.0001

# Row 141
$ws.Range("H141").Value = 10673.64
$ws.Range("I141").Value = 11656.637
$ws.Range("J141").Value = 3465
$ws.Range("K141").Value = 34969.911
$ws.Range("L141").Value = 10395
$ws.Range("M141").Value = -29789.911
$ws.Range("N141").Value = -20755

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2097.4119
$ws.Range("I74").Value = 1599.9656
$ws.Range("K74").Value = 1599.9656
$ws.Range("M74").Value = -725.9656

# Row 77
$ws.Range("H77").Value = 2097.4119
$ws.Range("I77").Value = 1599.9656
$ws.Range("K77").Value = 7999.828
$ws.Range("M77").Value = -3631.828

# Row 110
$ws.Range("H110").Value = 1744.3889
$ws.Range("I110").Value = 1812.9688
$ws.Range("J110").Value = 1195.75
$ws.Range("K110").Value = 1812.9688
$ws.Range("L110").Value = 1195.75
$ws.Range("M110").Value = 232.0311999999999
$ws.Range("N110").Value = -5285.75

# Row 132
$ws.Range("H132").Value = 4180.96
$ws.Range("I132").Value = 4153.2173
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 12459.6519
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -9929.651900000001
$ws.Range("N132").Value = -18560

# Row 141
$ws.Range("H141").Value = 95900
$ws.Range("J141").Value = 91800
$ws.Range("L141").Value = 91800
$ws.Range("N141").Value = -102160

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2148.8333
$ws.Range("I134").Value = 1823.6875
$ws.Range("K134").Value = 5471.0625
$ws.Range("M134").Value = -2936.0625

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5135.091
$ws.Range("I58").Value = 4602.7144
$ws.Range("J58").Value = 6066.75
$ws.Range("K58").Value = 4602.7144
$ws.Range("L58").Value = 6066.75
$ws.Range("M58").Value = -4399.7144
$ws.Range("N58").Value = -6472.75

# Row 132
$ws.Range("H132").Value = 11242.56
$ws.Range("I132").Value = 4311.5386
$ws.Range("K132").Value = 12934.6158
$ws.Range("M132").Value = -10404.6158

# Row 134
$ws.Range("H134").Value = 5292.125
$ws.Range("I134").Value = 5139.6113
$ws.Range("J134").Value = 5749.6665
$ws.Range("K134").Value = 15418.8339
$ws.Range("L134").Value = 17248.9995
$ws.Range("M134").Value = -12883.8339
$ws.Range("N134").Value = -22318.9995

# Row 136
$ws.Range("H136").Value = 5135.091
$ws.Range("I136").Value = 4602.7144
$ws.Range("J136").Value = 6066.75
$ws.Range("K136").Value = 13808.1432
$ws.Range("L136").Value = 18200.25
$ws.Range("M136").Value = -11258.1432
$ws.Range("N136").Value = -23300.25

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1177674.5
$ws.Range("I131").Value = 1401217.2
$ws.Range("K131").Value = 4203651.6
$ws.Range("M131").Value = -4198611.6

# Row 136
$ws.Range("H136").Value = 1461
$ws.Range("I136").Value = 1461
$ws.Range("K136").Value = 4383
$ws.Range("M136").Value = 717

# Row 139
$ws.Range("H139").Value = 2921.3809
$ws.Range("I139").Value = 1913.8334
$ws.Range("K139").Value = 5741.5002
$ws.Range("M139").Value = -601.5002000000004

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 38336000
$ws.Range("I80").Value = 230000000
$ws.Range("J80").Value = 3198.6
$ws.Range("K80").Value = 230000000
$ws.Range("L80").Value = 3198.6
$ws.Range("M80").Value = -229999002
$ws.Range("N80").Value = -5194.6

# Row 83
$ws.Range("H83").Value = 38336000
$ws.Range("I83").Value = 230000000
$ws.Range("J83").Value = 3198.6
$ws.Range("K83").Value = 1150000000
$ws.Range("L83").Value = 15993
$ws.Range("M83").Value = -1149995008
$ws.Range("N83").Value = -25977

# Row 102
$ws.Range("H102").Value = 4632
$ws.Range("I102").Value = 4285.207
$ws.Range("J102").Value = 5302.467
$ws.Range("K102").Value = 4285.207
$ws.Range("L102").Value = 5302.467
$ws.Range("M102").Value = -2663.207
$ws.Range("N102").Value = -8546.467000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 14998.333
$ws.Range("J43").Value = 14998.333
$ws.Range("L43").Value = 14998.333
$ws.Range("N43").Value = -15384.333

# Row 68
$ws.Range("H68").Value = 3146.7778
$ws.Range("I68").Value = 3146.7778
$ws.Range("K68").Value = 3146.7778
$ws.Range("M68").Value = -2397.7778

# Row 71
$ws.Range("H71").Value = 3146.7778
$ws.Range("I71").Value = 3146.7778
$ws.Range("K71").Value = 15733.889
$ws.Range("M71").Value = -11989.889

# Row 96
$ws.Range("H96").Value = 29999
$ws.Range("J96").Value = 29999
$ws.Range("L96").Value = 29999
$ws.Range("N96").Value = -35491

# Row 122
$ws.Range("H122").Value = 4271.136
$ws.Range("I122").Value = 4506.7856
$ws.Range("K122").Value = 13520.3568
$ws.Range("M122").Value = -11070.3568

# Row 132
$ws.Range("H132").Value = 19145.537
$ws.Range("I132").Value = 22686.922
$ws.Range("J132").Value = 9111.611000000001
$ws.Range("K132").Value = 68060.766
$ws.Range("L132").Value = 27334.833
$ws.Range("M132").Value = -65530.766
$ws.Range("N132").Value = -32394.833

# Row 136
$ws.Range("H136").Value = 7599.8667
$ws.Range("I136").Value = 9786.947
$ws.Range("J136").Value = 3822.182
$ws.Range("K136").Value = 29360.841
$ws.Range("L136").Value = 11466.546
$ws.Range("M136").Value = -26810.841
$ws.Range("N136").Value = -16566.546

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4568.2
$ws.Range("I132").Value = 4387.271
$ws.Range("K132").Value = 13161.813
$ws.Range("M132").Value = -10631.813

# Row 136
$ws.Range("H136").Value = 2820.4211
$ws.Range("I136").Value = 2769.8823
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 8309.6469
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -5759.6469
$ws.Range("N136").Value = -14850
